# Update the stats table on the "Axar Patel " sheet for Delhi Capitals.
# Cells C/D/E/F hold numbers-as-text (t="str" in the source OOXML), so we
# write each value with a leading apostrophe to force text entry (matching
# the existing "number stored as text" convention), then reset the cell
# style back to Normal so we don't leave a stray quote-prefix style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Row 2: C2 1->9, D2 4->9, E2 0->1 (F2 stays 0)
Set-TextValue "C2" "9"
Set-TextValue "D2" "9"
Set-TextValue "E2" "1"

# Row 4: C4 9->5, D4 7->6, F4 1->0 (E4 stays 0)
Set-TextValue "C4" "5"
Set-TextValue "D4" "6"
Set-TextValue "F4" "0"

# Row 5: C5 7->6, D5 4->9, E5 1->0 (F5 stays 0)
Set-TextValue "C5" "6"
Set-TextValue "D5" "9"
Set-TextValue "E5" "0"

# Row 7: C7 9->7, D7 9->4 (E7 stays 1)
Set-TextValue "C7" "7"
Set-TextValue "D7" "4"

# Row 8: C8 5->1, D8 6->4 (E8 stays 0)
Set-TextValue "C8" "1"
Set-TextValue "D8" "4"

# Row 9: C9 6->9, D9 9->7, F9 0->1 (E9 stays 0)
Set-TextValue "C9" "9"
Set-TextValue "D9" "7"
Set-TextValue "F9" "1"
